# Update new cluster checklist to reflect the deprecation of the SSO server.
#
# The old worksheet had a "SSO" section (server / app_id / app_secret /
# LDAP host / LDAP port / base query / uid field to query / username from
# field). This is replaced by a much smaller "Authentication" section
# (Authentication provider / credentials), and every section below it
# shifts up to follow directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Wipe out the whole old block (SSO .. ssh access .. Database password)
#    content AND formatting so we can rebuild it cleanly at the new row
#    numbers without leftover styling.
$ws.Range("A15:E44").Clear()

# 2. Authentication section (replaces SSO/LDAP section)
$ws.Range("A15").Value = "Authentication"
$ws.Range("A15").Font.Bold = $true

$ws.Range("A16").Value = "Authentication provider"
$ws.Range("B16").Value = "Google / OpenIDConnect / PAM / LDAP"

$ws.Range("A17").Value = "credentials"

# row 18 intentionally left blank

# 3. Keep storage (disk) section (shifted up from old row 26/27)
$ws.Range("A19").Value = "Keep storage (disk)"
$ws.Range("A19").Font.Bold = $true

$ws.Range("A20").Value = "hostname"
$ws.Range("A20").Font.Bold = $true
$ws.Range("B20").Value = "Volume partition/mount point"
$ws.Range("B20").Font.Bold = $true

# rows 21-22 intentionally left blank

# 4. slurm section (shifted up from old row 31-33)
$ws.Range("A24").Value = "slurm"
$ws.Range("A24").Font.Bold = $true

$ws.Range("A25").Value = "ControlMachine"
$ws.Range("A26").Value = "NodeName"

# row 27 intentionally left blank

# 5. ssh access section (shifted up from old row 36-39)
$ws.Range("A29").Value = "ssh access"
$ws.Range("A29").Font.Bold = $true

$ws.Range("A30").Value = "user"
$ws.Range("A31").Value = "key"
$ws.Range("A32").Value = "sudo password"

# 6. Database password (shifted up from old row 41)
$ws.Range("A34").Value = "Database password"
$ws.Range("A34").Font.Bold = $true

# rows 36-37 intentionally left blank (trailing spacer rows)

# 7. Move the active selection to reflect where the author was last
#    editing (A19, the new "Keep storage (disk)" heading).
$ws.Range("A19").Select()
